$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("education")

# Update the funding source for the Postdoctoral Research Stays grant
# from "Colciencias" to "Minciencias" (Colombia's science ministry was
# renamed from Colciencias to Minciencias).
$ws.Range("C2").Value = "Minciencias"

# Reflect the author's last selection position in the sheet view.
$ws.Range("C22").Select()
